$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# Insert a new column before column A, shifting everything right
$ws.Range("A1").EntireColumn.Insert()

# New SKU column header and values
$ws.Range("A1").Value = "sku"
$ws.Range("A2").Value = "SF-HEPA-H13-2424"
$ws.Range("A3").Value = "SF-PRE-G4-2020"

# Rename headers that changed wording (columns shifted right by one already)
$ws.Range("M1").Value = "estado"
$ws.Range("N1").Value = "etiquetas"
$ws.Range("O1").Value = "medida_nominal"
$ws.Range("P1").Value = "medida_real"
$ws.Range("V1").Value = "id_bind"
